$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): bold font size 18, centered horizontally and vertically ---
$ws.Range("A1:C1") | ForEach-Object {
    $_.HorizontalAlignment = -4108
    $_.VerticalAlignment = -4108
    $_.Font.Size = 18
}

# --- Data rows (2-37): Name/Local columns vertically centered, Descricao column wraps text ---
$ws.Range("A2:B37").VerticalAlignment = -4108
$ws.Range("C2:C37").WrapText = $true

# --- Fill in missing Descricao text for rows 25-37 ---

$ws.Range("C25").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   E aí, topa transformar ideias em potência com propósito de impactar a sociedade com a gente? Então faça parte de uma empresa que tem seus processos guiados por seus valores, valoriza seus colaboradores, sempre pensando no desenvolvimento profissional de cada um. E disponibiliza ambientes descontraídos, para que todos se sintam respeitados, acolhidos e acima de tudo: felizes.   Atualmente estamos com oportunidade para Desenvolvedor Full Stack Pleno Sênior   Descrição das atividades   O que você precisa ter: Experiência com .Net e Angular   Skill Técnico solicitado: Desejavel: Python, Xamarim, MongoDB, Postgre e Oracle   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C26").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   Atualmente estamos com oportunidade para Desenvolvedor FrontEnd React   Descrição das atividades   O que você precisa ter: - Sólidos conhecimentos técnicos em React - Sólidos conhecimentos técnicos em NextJS - Sólidos conhecimentos técnicos em Node - Experiência atuando com React em projetos grandes - Conhecimentos em GraphQL - Conhecimentos sobre arquitetura e boas práticas de desenvolvimento web - Noções sobre ambientes AWS   Skill Técnico solicitado: - Atuação em projetos de E-Commerce - Conhecimentos em React Native - Conhecimentos em Java (spring e springboot) - Conhecimentos em Kubernetes - Conhecimentos em PostgreSQL - Conhecimentos em Apache Cassandra - Conhecimentos/especializações em UI/UX   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C27").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   Atualmente estamos com oportunidade para Desenvolvedor Back End Java   Descrição das atividades   O que você precisa ter: Necessário: - Sólidos conhecimentos técnicos em Java 11 ao 14 - Sólidos conhecimentos técnicos em Spring, Springboot e microserviços - Experiência atuando com springboot e kubernetes em projetos grandes - Conhecimentos em Kubernetese Postgres - Conhecimentos sobre arquitetura e boas práticas de desenvolvimento backend - Noções sobre ambientes AWS   Skill Técnico solicitado: Interessante: - Atuação em projetos de E-Commerce - Certificações de Java - Conhecimentos em Apache Cassandra, React, React Native, Node, GraphQL Home Office – Tempo indeterminado   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C28").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   Atualmente estamos com oportunidade para Analista de Processos RPA   Descrição das atividades   O que você precisa ter: Mapear processos para automação juntamente com as áreas de negócio (AS-IS x TO-BE); Redigir documentação técnica de processos (PDD); Apoiar no pré venda de soluções RPA (reunião com clientes, elaboração de propostas, etc). Obrigatório: Experiência em mapeamento e modelagem de processos; Criação e manipulação de apresentações PowerPoint (.pptx) de qualidade visual; Conhecimentos em ferramentas de mapeamento de processos (ex: Bizagi, Visio, etc); Domínio da língua portuguesa e excelente comunicação oral e escrita. Desejável: Conhecimento em automação de processos RPA. Diferencial: Conhecimento em ferramentas de automação RPA (AutomationEdge, UiPath, etc Home Office – Tempo Indeterminado   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C29").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   E aí, topa transformar ideias em potência com propósito de impactar a sociedade com a gente? Então faça parte de uma empresa que tem seus processos guiados por seus valores, valoriza seus colaboradores, sempre pensando no desenvolvimento profissional de cada um. E disponibiliza ambientes descontraídos, para que todos se sintam respeitados, acolhidos e acima de tudo: felizes.   Atualmente estamos com oportunidade para Analista de NOC Jr III   Descrição das atividades   Skill Técnico solicitado: Ao mínimo 2 anos de experiência atuando no noc; Conhecimentos em Zabbix, Jira, Grafana e Solarwinds; Vivência em servidores Windows; Inglês Nível Intermediário. Diferencial: Conhecimentos em servidores Linux.   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C30").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   E aí, topa transformar ideias em potência com propósito de impactar a sociedade com a gente? Então faça parte de uma empresa que tem seus processos guiados por seus valores, valoriza seus colaboradores, sempre pensando no desenvolvimento profissional de cada um. E disponibiliza ambientes descontraídos, para que todos se sintam respeitados, acolhidos e acima de tudo: felizes.   Atualmente estamos com oportunidade para Dev Peoplesoft - ERP   Descrição das atividades   Skill Técnico solicitado: Application Designer, PeopleCode, Application Engine, Application Messaging, Component Interface, XML Publisher, SQR, Query, WorkFlow AWE, Component Interface, Application Package, Data Mover, File Layout, Programação PL/SQL. Conhecimento funcional desejável nos módulos do ERP Peoplesoft: o BI, AR, PO, AP, GL, TR, AM.   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C31").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   E aí, topa transformar ideias em potência com propósito de impactar a sociedade com a gente? Então faça parte de uma empresa que tem seus processos guiados por seus valores, valoriza seus colaboradores, sempre pensando no desenvolvimento profissional de cada um. E disponibiliza ambientes descontraídos, para que todos se sintam respeitados, acolhidos e acima de tudo: felizes.   Atualmente estamos com oportunidade para Estágio Logistica   Descrição das atividades   O que você terá que fazer: Controle de ativos que entram e saem do estoque.   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C32").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   E aí, topa transformar ideias em potência com propósito de impactar a sociedade com a gente? Então faça parte de uma empresa que tem seus processos guiados por seus valores, valoriza seus colaboradores, sempre pensando no desenvolvimento profissional de cada um. E disponibiliza ambientes descontraídos, para que todos se sintam respeitados, acolhidos e acima de tudo: felizes.   Atualmente estamos com oportunidade para Estágio TI   Descrição das atividades   O que você terá que fazer: "Profissional irá atuar na área de Suporte Técnico à Usuário N2, manutenção preventiva e corretiva de notebooks e periféricos. Realizará a instalação e configuração de aplicativo, verificação hardware e instalação de sistema operacional. Análises de equipamentos de troca ou desligamento, substituição de processadores, fontes, hd, memórias. Controle de ativos que entram e saem do estoque. • Área e especialização profissional: Informática, TI, Telecomunicações - TI, Logistica, Administração • Nível hierárquico: Estagiário • Local de trabalho: São Judas - SP   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida).'
$ws.Range("C33").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   Atualmente estamos com oportunidade para Especialista em QA   Descrição das atividades   O que você precisa ter: 5 anos + Experiencia em Desenvolvimento. Acompanhamento e realização de testes nas aplicações/softwares e aplicativos garantindo as qualidade e funcionalidades desenvolvidas com reportes das evoluções   O que você terá que fazer: . Realização dos testes progressivos e regressivos das versões dos aplicativos, backend e demais ferramentas de suporte. Especificação e automatização dos testes. Geração, integração e entrega das versões. Atuar na interface entre o cliente e a empresa na definição das metodologias e execução de testes. Definição e implementação de critérios de aceite e qualidade de software. Definição e implementação de métricas para medição dos processos de testes. Estudo de causa raiz, mapeamento dos incidentes problemas e endereçamento das correções. Execução e controle de Testes, validação das soluções e controle de versões   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C34").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   Atualmente estamos com oportunidade para Lider Técnico   Descrição das atividades   O que você terá que fazer: Interagir constantemente com os usuários de forma a garantir melhor compreensão; Desenvolver, testar e homologar solução definida; Identificar possíveis erros e demonstrar proatividade para resolvê-los da solução; Garantir em conjunto com a equipe o alcance das entregas definidas para a equipe; Liderar a equipe técnicamente, mostrando o caminho e resolvendo problemas; Requerimentos: Experiência com desenvolvimento de sistemas (Backend e Frontend) Experiência em .Net Core, Angular; Experiência com SQL Server; Vivência em times com desenvolvimento de metodologia Ágil e conhecimento contábil serão considerados diferenciais; Ter atuado como referente técnico ou líder de equipe técnica;   • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C35").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   Atualmente estamos com oportunidade para Analista Service Desk Jr I   Descrição das atividades   O que você terá que fazer: Receber os chamados de suporte abertos pelos usuários dos clientes cadastrados no sistema de gerenciamento ou encaminhados pelos Analistas de Service Desk; Realizar avaliação preliminar dos incidentes gerados, contatar o usuário para identificar a solução; Executar os processos para avaliação e identificação de problemas nos equipamentos dos usuários garantindo a sua satisfação;   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C36").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   Atualmente estamos com oportunidade para Desenvolvedor Back End Java – Springboot   Descrição das atividades   O que você terá que fazer: Atuar com a sustentação de APIs em Java 11 utilizando Springboot que estão publicadas na AWS como lambdas. - Desenvolver novas APIs em Java com springboot e publicá-las na AWS como lambdas.   Skill Técnico solicitado Necessário: - Experiência construindo e mantendo aplicações criadas em Java 8+ com Springboot. - Experiência utilizando o console da AWS. - Conhecimentos em OpenApi. - Conhecimentos em Lambdas (AWS). - Conhecimentos em API Gateway (AWS). Interessante: - Experiência com configuração de recursos da AWS. - Experiência criando ou mantendo aplicações Java com OpenApi. - Experiência com JSF. - Conhecimentos em Angular 6+. - Conhecimentos em testes unitários. Experiência: - Mínimo de 3 anos atuando na área de desenvolvimento. - Mínimo de 2 anos atuando com Springboot.   Dê um play em sua carreira! #VemSerCadmus.'
$ws.Range("C37").Value = 'Somos uma das melhores empresas para trabalhar pelo GPTW, por isso, buscamos os melhores talentos porque acreditamos que só com pessoas incríveis entregamos resultados incríveis.   E aí, topa transformar ideias em potência com propósito de impactar a sociedade com a gente? Então faça parte de uma empresa que tem seus processos guiados por seus valores, valoriza seus colaboradores, sempre pensando no desenvolvimento profissional de cada um. E disponibiliza ambientes descontraídos, para que todos se sintam respeitados, acolhidos e acima de tudo: felizes.   Atualmente estamos com oportunidade para Desenvolvedor FullStack   Nossos Diferencias: • Acompanhamento e desenvolvimento de carreira • Programa de pontos onde você escolhe seus benefícios (Plano de saúde, odontológico, vale alimentação, vale refeição, seguro de vida). • Treinamento e capacitação • Portal de conhecimentos e práticas • COP (Comunidade de Prática) • Clube de vantagens Lincard (Em casos de PJ ) • E muito mais….   Dê um play em sua carreira! #VemSerCadmus.'
